# Import ADJ column from excel budget file.
# Add three new header columns (ADJ1, ADj2, ADJ3) after the existing
# "MAR" column on the Budgets sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("S1").Value = "ADJ1"
$ws.Range("T1").Value = "ADj2"
$ws.Range("U1").Value = "ADJ3"

$ws.Range("U1").Select()
